$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row (row 28) - write in the order that reproduces the
# shared-string table order used by the target workbook:
# 26=/Tag, 27=kWh, 28=/h, 29=Brockdorf kW
$ws.Range("C28").Value = "/Tag"
$ws.Range("B28").Value = "kWh"
$ws.Range("D28").Value = "/h"
$ws.Range("E28").Value = "Brockdorf kW"

# New data row (row 29)
$ws.Range("B29").Value = 53000
$ws.Range("C29").Formula = "=B29/4"
$ws.Range("D29").Formula = "=C29/24"
$ws.Range("E29").Value = 1480000

# New result row (row 30) - written before the number format is applied to
# E29 so the format doesn't get inherited/auto-filled into E30
$ws.Range("E30").Formula = "=E29/D29"

# Apply the custom "#,##0" number format only to E29
$ws.Range("E29").NumberFormat = "#,##0"

# Move the active selection, matching the workbook as last edited
$ws.Range("E31").Select() | Out-Null
